# N5-kategori-teller-sorter.xlsx -- "Update 2020-04-23 part 1"
#
# Refresh the yearly counter tables on "Tellere kategorisert" (new 2006 row,
# corrected 2008-2011/2013 figures, two more rows for 2012/2015) and the
# category totals on "avskrivningsmaate" (new figures, "Besvart med brev"
# renamed to "Besvart med utgående brev", labels on sheet 1 getting a
# trailing colon).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tellere kategorisert")
$ws2 = $wb.Worksheets.Item("avskrivningsmaate")

# ---------------------------------------------------------------------
# Sheet "Tellere kategorisert"
# ---------------------------------------------------------------------

# Pre-fill rows 11:12 so the sort range below is non-empty, then run a
# (no-op content-wise) sort over A4:B12 purely so Excel widens the sheet's
# <sortState> bookkeeping to the new extent before we drop in the real
# figures (which are not in sorted order, so we set them explicitly after).
$ws1.Range("A11").Value = 0
$ws1.Range("B11").Value = 0
$ws1.Range("A12").Value = 0
$ws1.Range("B12").Value = 0

$sort1 = $ws1.Sort
$sort1.SortFields.Clear()
$sort1.SortFields.Add2($ws1.Range("A4"), 0, 1, $null, 0) | Out-Null
$sort1.SetRange($ws1.Range("A4:B12"))
$sort1.Header = 0
$sort1.Apply()

# Labels gain a trailing colon
$ws1.Range("A1").Value = "Sum hele:"
$ws1.Range("A2").Value = "Snitt hele:"
$ws1.Range("C1").Value = " ,Sum 2008-2011:"
$ws1.Range("C2").Value = " ,Snitt 2008-2011:"

# The 2008-2011 subtotal/sub-average now look at B6:B9 (rows shifted down
# by the new 2006 row at r4)
$ws1.Range("D1").Formula = "=SUM(B6:B9)"
$ws1.Range("D2").Formula = "=AVERAGE(B6:B8)"

# Year/count table, rows 4-12 (not in ascending year order -- matches the
# source data exactly)
$ws1.Range("A4").Value = 2006
$ws1.Range("B4").Value = 4
$ws1.Range("A5").Value = 2007
$ws1.Range("B5").Value = 50
$ws1.Range("A6").Value = 2008
$ws1.Range("B6").Value = 8264
$ws1.Range("A7").Value = 2009
$ws1.Range("B7").Value = 7911
$ws1.Range("A8").Value = 2011
$ws1.Range("B8").Value = 6205
$ws1.Range("A9").Value = 2010
$ws1.Range("B9").Value = 8727
$ws1.Range("A10").Value = 2013
$ws1.Range("B10").Value = 1
$ws1.Range("A11").Value = 2012
$ws1.Range("B11").Value = 2
$ws1.Range("A12").Value = 2015
$ws1.Range("B12").Value = 1

# ---------------------------------------------------------------------
# Sheet "avskrivningsmaate"
# ---------------------------------------------------------------------

# Same trick: run the sort (now against a custom list in real Excel) so
# <sortState> widens/re-keys to A1:B6, then overwrite with the real,
# updated figures/labels in their final (non-alphabetic) order.
$sort2 = $ws2.Sort
$sort2.SortFields.Clear()
$sort2.SortFields.Add2($ws2.Range("A1:A6"), 0, 1, $null, 0) | Out-Null
$sort2.SetRange($ws2.Range("A1:B6"))
$sort2.Header = 0
$sort2.Apply()

$ws2.Range("A1").Value = "Tatt til etterretning"
$ws2.Range("B1").Value = 8869
$ws2.Range("A2").Value = "Tatt til orientering"
$ws2.Range("B2").Value = 1178
$ws2.Range("A3").Value = "NULL"
$ws2.Range("B3").Value = 15
$ws2.Range("A4").Value = "Besvart med utgående brev"
$ws2.Range("B4").Value = 5661
$ws2.Range("A5").Value = "Besvart pr telefon"
$ws2.Range("B5").Value = 49
$ws2.Range("A6").Value = "Sak avsluttet"
$ws2.Range("B6").Value = 38

$wb.Application.Calculate()
